$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh: new Price (D) / Volume(1h) (E) figures
# for every listed coin. Rows 41-42 additionally swap places (MXToken
# now outranks RenderToken), so B/C (Coin name / Link) are rewritten too.
#
# Price cells that look like plain numbers ("0.9979", "313.27", ...) are
# pre-formatted as Text so Excel keeps them as literal strings (matching
# the source data, e.g. "1.200" / "0.01980" must keep their trailing
# zeros) instead of silently parsing them into numeric values.

$ws.Range("D2").Value = '27.403.51'
$ws.Range("E2").Value = '  +1.78%  '

$ws.Range("D3").Value = '1.828.33'
$ws.Range("E3").Value = '  +0.79%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9979'
$ws.Range("E4").Value = '  -0.38%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.27'
$ws.Range("E5").Value = '  +0.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9981'
$ws.Range("E6").Value = '  -0.42%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4474'
$ws.Range("E7").Value = '  -0.53%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3777'
$ws.Range("E8").Value = '  +2.59%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07409'
$ws.Range("E9").Value = '  +1.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8824'
$ws.Range("E10").Value = '  +3.39%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.92'
$ws.Range("E11").Value = '  +1.24%  '

$ws.Range("D12").Value = '1.823.64'
$ws.Range("E12").Value = '  +0.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.728'
$ws.Range("E13").Value = '  +1.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.435'
$ws.Range("E14").Value = '  +2.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.08'
$ws.Range("E15").Value = '  +1.29%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07074'
$ws.Range("E16").Value = '  -0.32%  '

$ws.Range("E17").Value = '  -0.30%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008826'
$ws.Range("E18").Value = '  +0.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9986'
$ws.Range("E19").Value = '  -0.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.07'
$ws.Range("E20").Value = '  +1.16%  '

$ws.Range("D21").Value = '27.376.26'
$ws.Range("E21").Value = '  +1.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.358'
$ws.Range("E22").Value = '  +4.02%  '

$ws.Range("E23").Value = '  +0.58%  '

$ws.Range("E24").Value = '  -1.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.89'
$ws.Range("E25").Value = '  -0.62%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.286'
$ws.Range("E26").Value = '  +2.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.67'
$ws.Range("E27").Value = '  +1.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.355'
$ws.Range("E28").Value = '  +2.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.38'
$ws.Range("E29").Value = '  +0.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08886'
$ws.Range("E30").Value = '  +0.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.7911'
$ws.Range("E31").Value = '  +5.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.200'
$ws.Range("E32").Value = '  +2.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.577'
$ws.Range("E33").Value = '  +3.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.918'
$ws.Range("E34").Value = '  -0.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9979'
$ws.Range("E35").Value = '  -0.43%  '

$ws.Range("E36").Value = '  +2.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01980'
$ws.Range("E37").Value = '  +1.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05272'
$ws.Range("E38").Value = '  +1.38%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.321'
$ws.Range("E39").Value = '  +3.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5306'
$ws.Range("E40").Value = '  +0.29%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.877'
$ws.Range("E41").Value = '  +0.35%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.335'
$ws.Range("E42").Value = '  +18.77%  '

$ws.Range("E43").Value = '  +0.71%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.664'
$ws.Range("E44").Value = '  +2.75%  '

$ws.Range("E45").Value = '  -3.15%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.57'
$ws.Range("E46").Value = '  +0.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '105.59'
$ws.Range("E47").Value = '  +0.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.687'
$ws.Range("E48").Value = '  +1.85%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9975'
$ws.Range("E49").Value = '  -0.45%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06396'
$ws.Range("E50").Value = '  +0.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '66.24'
$ws.Range("E51").Value = '  +5.80%  '
